# Automated tracker update: score a pending match and append the newest one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 55: the result for this pick is now known ---
$ws.Range("G55").Value = "Acierto"
$ws.Range("H55").Value = 3.33

# --- New row 72: latest match appended to the bottom of the tracker ---
$ws.Range("A72").Value = 14601440

# B72 is a literal date-like string ("2025-09-10"), not a real Excel date.
# Force Text format before assigning it so it isn't auto-converted into a
# date serial, then drop the temporary formatting again.
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = "2025-09-10"
$ws.Range("B72").ClearFormats()

$ws.Range("C72").Value = "Cezar Cretu"
$ws.Range("D72").Value = "Emilien Demanet"
$ws.Range("E72").Value = "Gana Emilien Demanet"
$ws.Range("F72").Value = 3

# G72/H72 (resultado/profit) are still unknown for the brand-new match, but
# the row keeps present-but-empty placeholder cells for those two columns,
# matching every other not-yet-settled row in the tracker. A lone "'" forces
# an explicit empty-text value instead of leaving the cell truly absent.
$ws.Range("G72").Value = "'"
$ws.Range("G72").ClearFormats()
$ws.Range("H72").Value = "'"
$ws.Range("H72").ClearFormats()
